$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text-looking numeric value as a shared string (no quotePrefix style),
# matching how the source workbook stores values like "1.7", "1.73", etc. as text.
function Set-TextValue($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Formula = "'" + $text
    $c.ClearFormats()
}

# Column J ("Media iterazioni") additions for the 12-case statistics update
Set-TextValue 68 10 "1.7"
Set-TextValue 70 10 "1.73"
Set-TextValue 72 10 "1.67"
Set-TextValue 75 10 "3.43"
Set-TextValue 77 10 "2.34"
Set-TextValue 79 10 "1.8"

$ws.Cells.Item(82, 10).Value = 1
$ws.Cells.Item(84, 10).Value = 1
$ws.Cells.Item(86, 10).Value = 1

# Update the view state (scrolled position / active selection) to match the edit session
$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Range("M88").Select() | Out-Null
